$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EngimonList")

# --- Update Battle Image (column F) for rows 7-20, in original authoring order ---
$ws.Range("F7").Value = "Engimon/077.png"
$ws.Range("F11").Value = "Engimon/184.png"
$ws.Range("F8").Value = "Engimon/148.png"
$ws.Range("F9").Value = "Engimon/350.png"
$ws.Range("F10").Value = "Engimon/446.png"
$ws.Range("F13").Value = "Engimon/485.png"
$ws.Range("F14").Value = "Engimon/479.png"
$ws.Range("F12").Value = "Engimon/392s.png"
$ws.Range("F15").Value = "Engimon/373.png"
$ws.Range("F17").Value = "Engimon/378.png"
$ws.Range("F16").Value = "Engimon/376.png"
$ws.Range("F18").Value = "Engimon/381.png"
$ws.Range("F19").Value = "Engimon/386.png"
$ws.Range("F20").Value = "Engimon/310.png"

# --- Update Sprite Image (column G) to the new character sprite; first write mints the shared string ---
$ws.Range("G2").Value = "Characters/boy_stand_south.png"

# --- Update Battle Image (column F) for rows 21-37, in original authoring order ---
$ws.Range("F21").Value = "Engimon/145.png"
$ws.Range("F22").Value = "Engimon/160.png"
$ws.Range("F23").Value = "Engimon/257.png"
$ws.Range("F24").Value = "Engimon/275fs.png"
$ws.Range("F25").Value = "Engimon/375.png"
$ws.Range("F26").Value = "Engimon/377.png"
$ws.Range("F27").Value = "Engimon/384.png"
$ws.Range("F28").Value = "Engimon/488.png"
$ws.Range("F29").Value = "Engimon/486.png"
$ws.Range("F30").Value = "Engimon/623.png"
$ws.Range("F31").Value = "Engimon/382.png"
$ws.Range("F32").Value = "Engimon/643.png"
$ws.Range("F36").Value = "Engimon/244.png"
$ws.Range("F34").Value = "Engimon/445.png"
$ws.Range("F33").Value = "Engimon/448s.png"
$ws.Range("F37").Value = "Engimon/065fs.png"

# --- Apply the same Sprite Image to the remaining rows ---
$ws.Range("G3").Value = "Characters/boy_stand_south.png"
$ws.Range("G4").Value = "Characters/boy_stand_south.png"
$ws.Range("G5").Value = "Characters/boy_stand_south.png"
$ws.Range("G6").Value = "Characters/boy_stand_south.png"
$ws.Range("G7").Value = "Characters/boy_stand_south.png"
$ws.Range("G8").Value = "Characters/boy_stand_south.png"
$ws.Range("G9").Value = "Characters/boy_stand_south.png"
$ws.Range("G10").Value = "Characters/boy_stand_south.png"
$ws.Range("G11").Value = "Characters/boy_stand_south.png"
$ws.Range("G12").Value = "Characters/boy_stand_south.png"
$ws.Range("G13").Value = "Characters/boy_stand_south.png"
$ws.Range("G14").Value = "Characters/boy_stand_south.png"
$ws.Range("G15").Value = "Characters/boy_stand_south.png"
$ws.Range("G16").Value = "Characters/boy_stand_south.png"
$ws.Range("G17").Value = "Characters/boy_stand_south.png"
$ws.Range("G18").Value = "Characters/boy_stand_south.png"
$ws.Range("G19").Value = "Characters/boy_stand_south.png"
$ws.Range("G20").Value = "Characters/boy_stand_south.png"
$ws.Range("G21").Value = "Characters/boy_stand_south.png"
$ws.Range("G22").Value = "Characters/boy_stand_south.png"
$ws.Range("G23").Value = "Characters/boy_stand_south.png"
$ws.Range("G24").Value = "Characters/boy_stand_south.png"
$ws.Range("G25").Value = "Characters/boy_stand_south.png"
$ws.Range("G26").Value = "Characters/boy_stand_south.png"
$ws.Range("G27").Value = "Characters/boy_stand_south.png"
$ws.Range("G28").Value = "Characters/boy_stand_south.png"
$ws.Range("G29").Value = "Characters/boy_stand_south.png"
$ws.Range("G30").Value = "Characters/boy_stand_south.png"
$ws.Range("G31").Value = "Characters/boy_stand_south.png"
$ws.Range("G32").Value = "Characters/boy_stand_south.png"
$ws.Range("G33").Value = "Characters/boy_stand_south.png"
$ws.Range("G34").Value = "Characters/boy_stand_south.png"
$ws.Range("G35").Value = "Characters/boy_stand_south.png"
$ws.Range("G36").Value = "Characters/boy_stand_south.png"
$ws.Range("G37").Value = "Characters/boy_stand_south.png"

# --- Restore view state: scroll so row 10 is at top, select F38 ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
[void]$ws.Range("F38").Select()
